# Regenerate save_data to use K (strikeouts) instead of Strike# for column G.
# The new K values are recalculated/re-scraped from source data and written
# directly into column G (rows 2-40), replacing the previous Strike# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 3
    5  = 0
    6  = 0
    7  = 2
    8  = 1
    9  = 0
    10 = 0
    11 = 2
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 1
    25 = 0
    26 = 2
    27 = 1
    28 = 2
    29 = 1
    30 = 2
    31 = 0
    32 = 2
    33 = 2
    34 = 2
    35 = 1
    36 = 0
    37 = 4
    38 = 3
    39 = 1
    40 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
